# Auto-generated Excel COM-interop script
# Refreshes market-price derived columns (H:N) for specific leve rows
# across all 8 job sheets, per the scheduled market-data runner update.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 351
$ws.Range("I31").Value = 351
$ws.Range("K31").Value = 1053
$ws.Range("M31").Value = -823
$ws.Range("H70").Value = 2694.5
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 2694.5
$ws.Range("K70").Value = 0
$ws.Range("L70").ClearContents()
$ws.Range("M70").Value = 8083.5
$ws.Range("N70").Value = -8623.5
$ws.Range("H73").Value = 2694.5
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 2694.5
$ws.Range("K73").Value = 0
$ws.Range("L73").ClearContents()
$ws.Range("M73").Value = 8083.5
$ws.Range("N73").Value = -9955.5
$ws.Range("H86").Value = 599.75
$ws.Range("I86").Value = 466.33334
$ws.Range("K86").Value = 466.33334
$ws.Range("M86").Value = 656.66666
$ws.Range("H89").Value = 599.75
$ws.Range("I89").Value = 466.33334
$ws.Range("K89").Value = 2331.6667
$ws.Range("M89").Value = 3284.3333
$ws.Range("H111").Value = 1214.5
$ws.Range("I111").Value = 1214.5
$ws.Range("K111").Value = 3643.5
$ws.Range("M111").Value = -576.5
$ws.Range("H125").Value = 17242
$ws.Range("I125").Value = 4100
$ws.Range("J125").Value = 25127.2
$ws.Range("K125").Value = 36900
$ws.Range("L125").Value = 226144.8
$ws.Range("M125").Value = -34440
$ws.Range("N125").Value = -231064.8
$ws.Range("H129").Value = 2942748.5
$ws.Range("I129").Value = 25000448
$ws.Range("J129").Value = 1722
$ws.Range("K129").Value = 75001344
$ws.Range("L129").Value = 5166
$ws.Range("M129").Value = -74996344
$ws.Range("N129").Value = -15166
$ws.Range("H135").Value = 850.5
$ws.Range("I135").Value = 850.5
$ws.Range("K135").Value = 7654.5
$ws.Range("M135").Value = -5119.5
$ws.Range("H138").Value = 4549038
$ws.Range("I138").Value = 25004624
$ws.Range("J138").Value = 3351.7222
$ws.Range("K138").Value = 75013872
$ws.Range("L138").Value = 10055.1666
$ws.Range("M138").Value = -75008732
$ws.Range("N138").Value = -20335.1666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9999.5
$ws.Range("I32").Value = 9999.5
$ws.Range("K32").Value = 9999.5
$ws.Range("M32").Value = -9712.5
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").ClearContents()
$ws.Range("N47").Value = 0
$ws.Range("H61").Value = 14665.667
$ws.Range("I61").Value = 10999.5
$ws.Range("J61").Value = 16498.75
$ws.Range("K61").Value = 10999.5
$ws.Range("L61").Value = 16498.75
$ws.Range("M61").Value = -10787.5
$ws.Range("N61").Value = -16922.75
$ws.Range("H74").Value = 4459.6665
$ws.Range("I74").Value = 1819.5
$ws.Range("K74").Value = 1819.5
$ws.Range("M74").Value = -945.5
$ws.Range("H77").Value = 4459.6665
$ws.Range("I77").Value = 1819.5
$ws.Range("K77").Value = 9097.5
$ws.Range("M77").Value = -4729.5
$ws.Range("H110").Value = 603.5
$ws.Range("I110").Value = 614.6667
$ws.Range("J110").Value = 570
$ws.Range("K110").Value = 614.6667
$ws.Range("L110").Value = 570
$ws.Range("M110").Value = 1430.3333
$ws.Range("N110").Value = -4660
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").ClearContents()
$ws.Range("N112").Value = 0
$ws.Range("H136").Value = 14665.667
$ws.Range("I136").Value = 10999.5
$ws.Range("J136").Value = 16498.75
$ws.Range("K136").Value = 32998.5
$ws.Range("L136").Value = 49496.25
$ws.Range("M136").Value = -30448.5
$ws.Range("N136").Value = -54596.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H29").Value = 3156.4
$ws.Range("I29").Value = 3156.4
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 3156.4
$ws.Range("L29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("N29").Value = -2867.4
$ws.Range("H37").Value = 2759.8
$ws.Range("I37").Value = 933
$ws.Range("K37").Value = 933
$ws.Range("M37").Value = -796
$ws.Range("H105").Value = 2247
$ws.Range("I105").Value = 1651
$ws.Range("K105").Value = 1651
$ws.Range("M105").Value = 96
$ws.Range("H110").Value = 12500
$ws.Range("J110").Value = 12500
$ws.Range("L110").Value = 12500
$ws.Range("N110").Value = -20680

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 249.92857
$ws.Range("I7").Value = 284.75
$ws.Range("J7").Value = 41
$ws.Range("K7").Value = 284.75
$ws.Range("L7").Value = 41
$ws.Range("M7").Value = -171.75
$ws.Range("N7").Value = -267
$ws.Range("H32").Value = 3505
$ws.Range("H41").Value = 59
$ws.Range("I41").Value = 59
$ws.Range("K41").Value = 59
$ws.Range("M41").Value = 369
$ws.Range("H60").Value = 18158.334
$ws.Range("J60").Value = 25000
$ws.Range("L60").Value = 25000
$ws.Range("N60").Value = -26022
$ws.Range("H105").Value = 1503.3334
$ws.Range("J105").Value = 500
$ws.Range("L105").Value = 500
$ws.Range("N105").Value = -3994
$ws.Range("H132").Value = 4325.5
$ws.Range("I132").Value = 2838.5454
$ws.Range("K132").Value = 8515.636200000001
$ws.Range("M132").Value = -5985.636200000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("L17").ClearContents()
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = 0
$ws.Range("H69").Value = 512
$ws.Range("I69").Value = 512
$ws.Range("K69").Value = 1536
$ws.Range("M69").Value = -725
$ws.Range("H72").Value = 512
$ws.Range("I72").Value = 512
$ws.Range("K72").Value = 4608
$ws.Range("M72").Value = -552
$ws.Range("H137").Value = 1999
$ws.Range("J137").Value = 1999
$ws.Range("L137").Value = 5997
$ws.Range("N137").Value = -16197
$ws.Range("H139").Value = 2959.4
$ws.Range("I139").Value = 2959.4
$ws.Range("K139").Value = 8878.200000000001
$ws.Range("M139").Value = -3738.200000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 1004
$ws.Range("I5").Value = 1004
$ws.Range("K5").Value = 1004
$ws.Range("M5").Value = -892

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H30").Value = 1430
$ws.Range("I30").Value = 662.5
$ws.Range("J30").Value = 4500
$ws.Range("K30").Value = 662.5
$ws.Range("L30").Value = 4500
$ws.Range("M30").Value = -554.5
$ws.Range("N30").Value = -4716
$ws.Range("H46").Value = 6319.4
$ws.Range("I46").Value = 5566.6665
$ws.Range("J46").Value = 7448.5
$ws.Range("K46").Value = 5566.6665
$ws.Range("L46").Value = 7448.5
$ws.Range("M46").Value = -5378.6665
$ws.Range("N46").Value = -7824.5
$ws.Range("H80").Value = 32500
$ws.Range("J80").Value = 32500
$ws.Range("L80").Value = 32500
$ws.Range("N80").Value = -34746
$ws.Range("H83").Value = 32500
$ws.Range("J83").Value = 32500
$ws.Range("L83").Value = 97500
$ws.Range("N83").Value = -108732
$ws.Range("H122").Value = 1499.3334
$ws.Range("I122").Value = 1249.5
$ws.Range("J122").Value = 1999
$ws.Range("K122").Value = 3748.5
$ws.Range("L122").Value = 5997
$ws.Range("M122").Value = -1298.5
$ws.Range("N122").Value = -10897
$ws.Range("H132").Value = 12494.6
$ws.Range("I132").Value = 12494.6
$ws.Range("K132").Value = 37483.8
$ws.Range("M132").Value = -34953.8
$ws.Range("H136").Value = 12928.286
$ws.Range("I136").Value = 7624.75
$ws.Range("K136").Value = 22874.25
$ws.Range("M136").Value = -20324.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 39500
$ws.Range("I2").Value = 29000
$ws.Range("K2").Value = 29000
$ws.Range("M2").Value = -28888
$ws.Range("H81").Value = 492.5
$ws.Range("I81").Value = 492.5
$ws.Range("K81").Value = 985
$ws.Range("M81").Value = 76
$ws.Range("H84").Value = 492.5
$ws.Range("I84").Value = 492.5
$ws.Range("K84").Value = 4925
$ws.Range("M84").Value = 379
$ws.Range("H100").Value = 225
$ws.Range("I100").Value = 225
$ws.Range("K100").Value = 450
$ws.Range("M100").Value = 91
$ws.Range("H103").Value = 14000.333
$ws.Range("J103").Value = 14000.333
$ws.Range("L103").Value = 14000.333
$ws.Range("N103").Value = -16344.333
$ws.Range("H132").Value = 1249.6666
$ws.Range("I132").Value = 1249.6666
$ws.Range("K132").Value = 3748.9998
$ws.Range("M132").Value = -1218.9998

